$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-29 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-30 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("39×79=3081", $true, $false, $false, $false, $false, $true, 1, $false, "32×34=1088", 2) | Out-Null
$d.Content.Find.Execute("14×38=532", $true, $false, $false, $false, $false, $true, 1, $false, "70×98=6860", 2) | Out-Null
$d.Content.Find.Execute("59×35=2065", $true, $false, $false, $false, $false, $true, 1, $false, "24×19=456", 2) | Out-Null
$d.Content.Find.Execute("79×91=7189", $true, $false, $false, $false, $false, $true, 1, $false, "14×87=1218", 2) | Out-Null
$d.Content.Find.Execute("60×92=5520", $true, $false, $false, $false, $false, $true, 1, $false, "78×54=4212", 2) | Out-Null
$d.Content.Find.Execute("56×98=5488", $true, $false, $false, $false, $false, $true, 1, $false, "86×86=7396", 2) | Out-Null
$d.Content.Find.Execute("74×43=3182", $true, $false, $false, $false, $false, $true, 1, $false, "11×44=484", 2) | Out-Null
$d.Content.Find.Execute("78×93=7254", $true, $false, $false, $false, $false, $true, 1, $false, "70×71=4970", 2) | Out-Null
$d.Content.Find.Execute("57×58=3306", $true, $false, $false, $false, $false, $true, 1, $false, "67×79=5293", 2) | Out-Null
$d.Content.Find.Execute("74×97=7178", $true, $false, $false, $false, $false, $true, 1, $false, "30×69=2070", 2) | Out-Null
$d.Content.Find.Execute("66×69=4554", $true, $false, $false, $false, $false, $true, 1, $false, "28×89=2492", 2) | Out-Null
$d.Content.Find.Execute("57×41=2337", $true, $false, $false, $false, $false, $true, 1, $false, "29×48=1392", 2) | Out-Null
$d.Content.Find.Execute("80×75=6000", $true, $false, $false, $false, $false, $true, 1, $false, "99×84=8316", 2) | Out-Null
$d.Content.Find.Execute("13×44=572", $true, $false, $false, $false, $false, $true, 1, $false, "81×71=5751", 2) | Out-Null
$d.Content.Find.Execute("82×41=3362", $true, $false, $false, $false, $false, $true, 1, $false, "47×18=846", 2) | Out-Null
$d.Content.Find.Execute("88×46=4048", $true, $false, $false, $false, $false, $true, 1, $false, "58×84=4872", 2) | Out-Null
$d.Content.Find.Execute("43×87=3741", $true, $false, $false, $false, $false, $true, 1, $false, "17×78=1326", 2) | Out-Null
$d.Content.Find.Execute("69×87=6003", $true, $false, $false, $false, $false, $true, 1, $false, "21×54=1134", 2) | Out-Null
$d.Content.Find.Execute("55×94=5170", $true, $false, $false, $false, $false, $true, 1, $false, "91×74=6734", 2) | Out-Null
$d.Content.Find.Execute("33×30=990", $true, $false, $false, $false, $false, $true, 1, $false, "42×34=1428", 2) | Out-Null
$d.Content.Find.Execute("51×81=4131", $true, $false, $false, $false, $false, $true, 1, $false, "66×98=6468", 2) | Out-Null
$d.Content.Find.Execute("84×11=924", $true, $false, $false, $false, $false, $true, 1, $false, "48×64=3072", 2) | Out-Null
$d.Content.Find.Execute("17×45=765", $true, $false, $false, $false, $false, $true, 1, $false, "15×62=930", 2) | Out-Null
$d.Content.Find.Execute("37×17=629", $true, $false, $false, $false, $false, $true, 1, $false, "65×70=4550", 2) | Out-Null
$d.Content.Find.Execute("42×70=2940", $true, $false, $false, $false, $false, $true, 1, $false, "64×64=4096", 2) | Out-Null
